# Generate Report for Handoff
#
# Swaps the "59b68d80..." and "aca42e72..." rows (rows 2/3) across the
# Overview / zh-cn / de-de sheets so that "aca42e72..." now occupies row 2
# and "59b68d80..." occupies row 3, and refreshes the "59b68d80..." entry's
# status/timestamps/error-detail to reflect that it is ready for a fresh
# handoff (no longer "in sync", but out of date vs. the latest source).

function Set-HyperlinkDisplay {
    param($ws, $addr, $text)
    foreach ($h in $ws.Hyperlinks) {
        if ($h.Range.Address() -eq $addr) {
            $h.TextToDisplay = $text
        }
    }
}

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------
$ov = $wb.Worksheets.Item("Overview")

$ov.Range("A2").Value = "aca42e72-7742-48c2-946a-9507802b42c7.md"
$ov.Range("B2").Value = "e2e\aca42e72-7742-48c2-946a-9507802b42c7.md"

$ov.Range("A3").Value = "59b68d80-0cb3-4ce2-9165-59d760900154.md"
$ov.Range("B3").Value = "e2e\59b68d80-0cb3-4ce2-9165-59d760900154.md"
$ov.Range("E3").Value = "Ready for handoff"
$ov.Range("F3").Value = "Ready for handoff"
$ov.Range("G3").Value = "2016-09-06 11:04:27"

Set-HyperlinkDisplay $ov '$B$2' "e2e\aca42e72-7742-48c2-946a-9507802b42c7.md"
Set-HyperlinkDisplay $ov '$B$3' "e2e\59b68d80-0cb3-4ce2-9165-59d760900154.md"

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$zh = $wb.Worksheets.Item("zh-cn")

$zh.Range("A2").Value = "aca42e72-7742-48c2-946a-9507802b42c7.md"
$zh.Range("G2").Value = "aca42e72-7742-48c2-946a-9507802b42c7.14d05f70bdbaf7ef0d9195b1fff84847f8843526.zh-cn.xlf"
$zh.Range("I2").Value = "aca42e72-7742-48c2-946a-9507802b42c7.md"
$zh.Range("J2").Value = "aca42e72-7742-48c2-946a-9507802b42c7.14d05f70bdbaf7ef0d9195b1fff84847f8843526.zh-cn.xlf"

$zh.Range("A3").Value = "59b68d80-0cb3-4ce2-9165-59d760900154.md"
$zh.Range("C3").Value = "Ready for handoff"
$zh.Range("G3").Value = "59b68d80-0cb3-4ce2-9165-59d760900154.0f271ebaa6a94107966737bf332c97a969b0b72b.zh-cn.xlf"
$zh.Range("H3").Value = "2016-09-06 11:04:23"
$zh.Range("I3").Value = "59b68d80-0cb3-4ce2-9165-59d760900154.md"
$zh.Range("J3").Value = "59b68d80-0cb3-4ce2-9165-59d760900154.0f271ebaa6a94107966737bf332c97a969b0b72b.zh-cn.xlf"
$zh.Range("P3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/6b5374ccd07cbd13644b31b4872ddc37a61171b5/e2e/59b68d80-0cb3-4ce2-9165-59d760900154.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/c6d954f7b508f1107f1e1bd490a94d29fc389033/e2e/59b68d80-0cb3-4ce2-9165-59d760900154.md."

Set-HyperlinkDisplay $zh '$A$2' "aca42e72-7742-48c2-946a-9507802b42c7.md"
Set-HyperlinkDisplay $zh '$I$2' "aca42e72-7742-48c2-946a-9507802b42c7.md"
Set-HyperlinkDisplay $zh '$A$3' "59b68d80-0cb3-4ce2-9165-59d760900154.md"
Set-HyperlinkDisplay $zh '$I$3' "59b68d80-0cb3-4ce2-9165-59d760900154.md"

$zh.Columns.Item(16).ColumnWidth = 39.17

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$de = $wb.Worksheets.Item("de-de")

$de.Range("A2").Value = "aca42e72-7742-48c2-946a-9507802b42c7.md"
$de.Range("G2").Value = "aca42e72-7742-48c2-946a-9507802b42c7.14d05f70bdbaf7ef0d9195b1fff84847f8843526.de-de.xlf"
$de.Range("I2").Value = "aca42e72-7742-48c2-946a-9507802b42c7.md"
$de.Range("J2").Value = "aca42e72-7742-48c2-946a-9507802b42c7.14d05f70bdbaf7ef0d9195b1fff84847f8843526.de-de.xlf"

$de.Range("A3").Value = "59b68d80-0cb3-4ce2-9165-59d760900154.md"
$de.Range("C3").Value = "Ready for handoff"
$de.Range("G3").Value = "59b68d80-0cb3-4ce2-9165-59d760900154.0f271ebaa6a94107966737bf332c97a969b0b72b.de-de.xlf"
$de.Range("H3").Value = "2016-09-06 11:04:27"
$de.Range("I3").Value = "59b68d80-0cb3-4ce2-9165-59d760900154.md"
$de.Range("J3").Value = "59b68d80-0cb3-4ce2-9165-59d760900154.0f271ebaa6a94107966737bf332c97a969b0b72b.de-de.xlf"
$de.Range("P3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/6b5374ccd07cbd13644b31b4872ddc37a61171b5/e2e/59b68d80-0cb3-4ce2-9165-59d760900154.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/c6d954f7b508f1107f1e1bd490a94d29fc389033/e2e/59b68d80-0cb3-4ce2-9165-59d760900154.md."

Set-HyperlinkDisplay $de '$A$2' "aca42e72-7742-48c2-946a-9507802b42c7.md"
Set-HyperlinkDisplay $de '$I$2' "aca42e72-7742-48c2-946a-9507802b42c7.md"
Set-HyperlinkDisplay $de '$A$3' "59b68d80-0cb3-4ce2-9165-59d760900154.md"
Set-HyperlinkDisplay $de '$I$3' "59b68d80-0cb3-4ce2-9165-59d760900154.md"

$de.Columns.Item(16).ColumnWidth = 39.17
